$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.370.61'
$ws.Range("E2").Value = '  +3.59%  '
$ws.Range("D3").Value = '1.836.92'
$ws.Range("E3").Value = '  +3.64%  '
$ws.Range("E4").Value = '  +2.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '316.85'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.82%  '
$ws.Range("E6").Value = '  +1.74%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4355'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.99%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3720'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.69%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07332'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.33%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8717'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.77%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '21.33'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +4.54%  '
$ws.Range("D12").Value = '1.939.32'
$ws.Range("E12").Value = '  +7.72%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.465'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +4.13%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.680'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.72%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.07124'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.23%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '82.10'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.20%  '
$ws.Range("E17").Value = '  +2.02%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008967'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.29%  '
$ws.Range("E19").Value = '  +1.89%  '
$ws.Range("E20").Value = '  +2.85%  '
$ws.Range("D21").Value = '27.401.05'
$ws.Range("E21").Value = '  +3.60%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.244'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.79%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.13'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.66%  '
$ws.Range("D24").Value = '2.138.79'
$ws.Range("E24").Value = '  +6.04%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '156.64'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.91%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.893'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +4.62%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.53'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.92%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.236'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.34%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.920'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +8.57%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '115.44'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.60%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09037'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.65%  '
$ws.Range("E32").Value = '  +7.96%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7587'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +4.76%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.456'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.18%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.858'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.93%  '
$ws.Range("E36").Value = '  +2.01%  '
$ws.Range("E37").Value = '  +4.43%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01956'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.83%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05238'
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5161'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +5.00%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.778'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +7.52%  '
$ws.Range("E42").Value = '  +2.89%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.531'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.94%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.453'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +6.21%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '108.33'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.48%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.47'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.65%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.025'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.08%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.671'
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.4620'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +4.20%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06293'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.86%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.879'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +9.60%  '
